# Added room assignment between rooms and tenants
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing row 8 (tenant G) gains a tenant_friends value
$ws.Range("D8").Value = "A"

# New tenant/room-assignment rows 9-15
$data = @(
    @("H", "H", 2, "G"),
    @("I", "I", 2, "J, J, L"),
    @("J", "J", 2, "L, M, G"),
    @("K", "K", 1, "A, B, C"),
    @("L", "L", 2, "D, G"),
    @("M", "M", 1, "N"),
    @("N", "N", 1, $null)
)

$r = 9
foreach ($row in $data) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    if ($row[3] -ne $null) {
        $ws.Range("D$r").Value = $row[3]
    }
    $r = $r + 1
}

# Copy the formatting from the already-formatted row (row 8's A:C, which
# still carries the original style) down to every new/changed cell so they
# all pick up the same cell style as the rest of the table
# (xlPasteFormats = -4122). Row 15 has no tenant_friends entry, so its
# D cell must be left untouched (no cell should be created there).
$ws.Range("C8").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null
$ws.Range("A8:D8").Copy() | Out-Null
$ws.Range("A9:D14").PasteSpecial(-4122) | Out-Null
$ws.Range("A8:C8").Copy() | Out-Null
$ws.Range("A15:C15").PasteSpecial(-4122) | Out-Null
